# chore: update Sheets via scheduled runner
# Refresh the pricing/profit columns (H:N) on the Ultima_Profits workbook
# for the rows whose source market data changed in this run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3620.9375
$ws.Range("I51").Value = 3200
$ws.Range("J51").Value = 3718.077
$ws.Range("K51").Value = 3200
$ws.Range("L51").Value = 3718.077
$ws.Range("M51").Value = -2716
$ws.Range("N51").Value = -4686.077

$ws.Range("H116").Value = 2237.4375
$ws.Range("I116").Value = 1856
$ws.Range("J116").Value = 2727.8572
$ws.Range("K116").Value = 1856
$ws.Range("L116").Value = 2727.8572
$ws.Range("M116").Value = 1586
$ws.Range("N116").Value = -9611.8572

$ws.Range("H132").Value = 4402.4365
$ws.Range("I132").Value = 3802.68
$ws.Range("J132").Value = 10400
$ws.Range("K132").Value = 11408.04
$ws.Range("L132").Value = 31200
$ws.Range("M132").Value = -8878.039999999999
$ws.Range("N132").Value = -36260

$ws.Range("H138").Value = 2075.5059
$ws.Range("I138").Value = 1210.8723
$ws.Range("J138").Value = 3091.45
$ws.Range("K138").Value = 3632.6169
$ws.Range("L138").Value = 9274.349999999999
$ws.Range("M138").Value = 1507.3831
$ws.Range("N138").Value = -19554.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 636.8627300000001
$ws.Range("I2").Value = 558.9048
$ws.Range("J2").Value = 1000.6667
$ws.Range("K2").Value = 558.9048
$ws.Range("L2").Value = 1000.6667
$ws.Range("M2").Value = -445.9048
$ws.Range("N2").Value = -1226.6667

$ws.Range("H74").Value = 15626922
$ws.Range("I74").Value = 27778928
$ws.Range("J74").Value = 2914.5715
$ws.Range("K74").Value = 27778928
$ws.Range("L74").Value = 2914.5715
$ws.Range("M74").Value = -27778054
$ws.Range("N74").Value = -4662.5715

$ws.Range("H77").Value = 15626922
$ws.Range("I77").Value = 27778928
$ws.Range("J77").Value = 2914.5715
$ws.Range("K77").Value = 138894640
$ws.Range("L77").Value = 14572.8575
$ws.Range("M77").Value = -138890272
$ws.Range("N77").Value = -23308.8575

$ws.Range("H116").Value = 636.8627300000001
$ws.Range("I116").Value = 558.9048
$ws.Range("J116").Value = 1000.6667
$ws.Range("K116").Value = 558.9048
$ws.Range("L116").Value = 1000.6667
$ws.Range("M116").Value = 1735.0952
$ws.Range("N116").Value = -5588.6667

$ws.Range("H132").Value = 4547130.5
$ws.Range("I132").Value = 5815439.5
$ws.Range("K132").Value = 17446318.5
$ws.Range("M132").Value = -17443788.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 636.8627300000001
$ws.Range("I3").Value = 558.9048
$ws.Range("J3").Value = 1000.6667
$ws.Range("K3").Value = 558.9048
$ws.Range("L3").Value = 1000.6667
$ws.Range("M3").Value = -444.9048
$ws.Range("N3").Value = -1228.6667

$ws.Range("H69").Value = 43628.332
$ws.Range("J69").Value = 43628.332
$ws.Range("L69").Value = 43628.332
$ws.Range("N69").Value = -45250.332

$ws.Range("H72").Value = 43628.332
$ws.Range("J72").Value = 43628.332
$ws.Range("L72").Value = 130884.996
$ws.Range("N72").Value = -138996.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1772.9375
$ws.Range("I105").Value = 1494.7
$ws.Range("J105").Value = 2236.6667
$ws.Range("K105").Value = 1494.7
$ws.Range("L105").Value = 2236.6667
$ws.Range("M105").Value = 252.3
$ws.Range("N105").Value = -5730.6667

$ws.Range("H132").Value = 11365812
$ws.Range("J132").Value = 3861.1
$ws.Range("L132").Value = 11583.3
$ws.Range("N132").Value = -16643.3

$ws.Range("H134").Value = 1704902.9
$ws.Range("I134").Value = 7227.2
$ws.Range("J134").Value = 2648056
$ws.Range("K134").Value = 21681.6
$ws.Range("L134").Value = 7944168
$ws.Range("M134").Value = -19146.6
$ws.Range("N134").Value = -7949238

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H105").Value = 4250
$ws.Range("J105").Value = 4250
$ws.Range("L105").Value = 12750
$ws.Range("N105").Value = -17992

$ws.Range("H121").Value = 589.9231
$ws.Range("I121").Value = 433.75
$ws.Range("J121").Value = 839.8
$ws.Range("K121").Value = 1301.25
$ws.Range("L121").Value = 2519.4
$ws.Range("M121").Value = 8.75
$ws.Range("N121").Value = -5139.4

$ws.Range("H123").Value = 4689.92
$ws.Range("I123").Value = 1883.1818
$ws.Range("J123").Value = 6895.2144
$ws.Range("K123").Value = 5649.5454
$ws.Range("L123").Value = 20685.6432
$ws.Range("M123").Value = -3199.5454
$ws.Range("N123").Value = -25585.6432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8845.25
$ws.Range("I70").Value = 10255.549
$ws.Range("J70").Value = 3987.5557
$ws.Range("K70").Value = 10255.549
$ws.Range("L70").Value = 3987.5557
$ws.Range("M70").Value = -9985.549000000001
$ws.Range("N70").Value = -4527.5557

$ws.Range("H73").Value = 8845.25
$ws.Range("I73").Value = 10255.549
$ws.Range("J73").Value = 3987.5557
$ws.Range("K73").Value = 10255.549
$ws.Range("L73").Value = 3987.5557
$ws.Range("M73").Value = -9319.549000000001
$ws.Range("N73").Value = -5859.5557

$ws.Range("H97").Value = 1533.875
$ws.Range("I97").Value = 1850
$ws.Range("J97").Value = 1344.2
$ws.Range("K97").Value = 1850
$ws.Range("L97").Value = 1344.2
$ws.Range("M97").Value = -1354
$ws.Range("N97").Value = -2336.2

$ws.Range("H122").Value = 3336165.2
$ws.Range("I122").Value = 6062701.5
$ws.Range("J122").Value = 3732
$ws.Range("K122").Value = 18188104.5
$ws.Range("L122").Value = 11196
$ws.Range("M122").Value = -18185654.5
$ws.Range("N122").Value = -16096

$ws.Range("H132").Value = 4385.548
$ws.Range("I132").Value = 4717.1763
$ws.Range("J132").Value = 2976.125
$ws.Range("K132").Value = 14151.5289
$ws.Range("L132").Value = 8928.375
$ws.Range("M132").Value = -11621.5289
$ws.Range("N132").Value = -13988.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5446.8667
$ws.Range("I62").Value = 5242.857
$ws.Range("J62").Value = 5625.375
$ws.Range("K62").Value = 5242.857
$ws.Range("L62").Value = 5625.375
$ws.Range("M62").Value = -4618.857
$ws.Range("N62").Value = -6873.375

$ws.Range("H65").Value = 5446.8667
$ws.Range("I65").Value = 5242.857
$ws.Range("J65").Value = 5625.375
$ws.Range("K65").Value = 26214.285
$ws.Range("L65").Value = 28126.875
$ws.Range("M65").Value = -23094.285
$ws.Range("N65").Value = -34366.875

$ws.Range("H107").Value = 970.64703
$ws.Range("I107").Value = 1343.5555
$ws.Range("J107").Value = 551.125
$ws.Range("K107").Value = 4030.6665
$ws.Range("L107").Value = 1653.375
$ws.Range("M107").Value = -2110.6665
$ws.Range("N107").Value = -5493.375

$ws.Range("H136").Value = 1008.53125
$ws.Range("I136").Value = 1009.7586
$ws.Range("K136").Value = 3029.2758
$ws.Range("M136").Value = -479.2757999999999

Write-Host "Updated profit columns across ALC, ARM, BSM, CRP, CUL, GSM, WVR sheets"
